# Complementação da lista com brainstorm
# Adds one new "brainstormed" idea to the requirements list (row 31,
# leaving row 30 blank like the existing gap after row 2/3) and marks it
# in red to flag it as a new/unvetted suggestion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A31")
$cell.Value = ".Integração com sistema da polícia sobre carros roubados"
$cell.Font.Color = 192

# Keep the selection near the bottom of the list, matching where the
# author's cursor ended up after adding the new item.
$ws.Range("A29").Select()
